$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header cells: "_old" -> "_FV2304", "_new" -> "_FV2310"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Text
    if ($val.EndsWith("_old")) {
        $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2304"
    } elseif ($val.EndsWith("_new")) {
        $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2310"
    }
}

# 2) Add table (ListObject) over A1:U58 with headers
$range = $ws.Range("A1:U58")
$table = $ws.ListObjects.Add(1, $range, [System.Reflection.Missing]::Value, 1)
$table.Name = "Table1"

# 3) Freeze header row (pane split) and set selection on the frozen pane
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
